$wb = $excel.ActiveWorkbook

# Mapping of row -> new F value, to be applied on both "展览" and "全部类型" sheets
$updates = @{
    3  = 1345
    4  = 151
    11 = 4541
    12 = 6805
    16 = 569
    18 = 4124
    19 = 491
    20 = 72
    24 = 548
    25 = 168
    27 = 360
    29 = 224
    36 = 544
    40 = 73
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
